$wb = $excel.ActiveWorkbook

# Work on the "Repayment schedule" sheet
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N, shifting N:P -> O:Q
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = 10.17

# Select a cell on this sheet and make it the active sheet/tab
$ws.Activate()
$ws.Range("R7").Select()

# The "NewLoanInput" sheet should no longer be the tab-selected sheet
